# Apply the cryptos-list price/volume refresh described by the commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking strings (e.g. thousands-dot prices) that must
# stay text, exactly as authored, so every Price write is apostrophe-prefixed
# (the same trick Excel itself uses to force text entry).

$ws.Range('D2').Value = '''26.761.13'
$ws.Range('E2').Value = '  -2.52%  '
$ws.Range('D3').Value = '''1.566.69'
$ws.Range('E3').Value = '  +0.12%  '
$ws.Range('E4').Value = '  +0.08%  '
$ws.Range('D5').Value = '''206.43'
$ws.Range('E5').Value = '  -0.94%  '
$ws.Range('D6').Value = '''0.489'
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('D8').Value = '''21.90'
$ws.Range('E8').Value = '  -0.83%  '
$ws.Range('D9').Value = '''0.248'
$ws.Range('E9').Value = '  -0.58%  '
$ws.Range('E10').Value = '  -1.26%  '
$ws.Range('E11').Value = '  -0.18%  '
$ws.Range('D12').Value = '''1.788.39'
$ws.Range('E12').Value = '  +0.20%  '
$ws.Range('D13').Value = '''1.573.97'
$ws.Range('E13').Value = '  +0.60%  '
$ws.Range('E14').Value = '  -2.39%  '
$ws.Range('D15').Value = '''0.516'
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('D16').Value = '''61.55'
$ws.Range('E16').Value = '  -3.14%  '
$ws.Range('D17').Value = '''26.787.65'
$ws.Range('D18').Value = '''214.80'
$ws.Range('E18').Value = '  +0.88%  '
$ws.Range('D19').Value = '''7.37'
$ws.Range('E19').Value = '  +1.73%  '
$ws.Range('E20').Value = '  -1.75%  '
$ws.Range('E21').Value = '  +0.14%  '
$ws.Range('E22').Value = '  -0.06%  '
$ws.Range('E23').Value = '  -2.75%  '
$ws.Range('E24').Value = '  -1.59%  '
$ws.Range('D25').Value = '''152.62'
$ws.Range('E25').Value = '  -0.13%  '
$ws.Range('E26').Value = '  +1.07%  '
$ws.Range('D27').Value = '''14.92'
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('E29').Value = '  -1.27%  '
$ws.Range('E30').Value = '  -1.21%  '
$ws.Range('E31').Value = '  -3.69%  '
$ws.Range('E32').Value = '  -1.35%  '
$ws.Range('D33').Value = '''1.388.38'
$ws.Range('E33').Value = '  +1.06%  '
$ws.Range('E34').Value = '  -1.13%  '
$ws.Range('E35').Value = '  +0.67%  '
$ws.Range('E36').Value = '  -0.89%  '
$ws.Range('D37').Value = '''0.924'
$ws.Range('E37').Value = '  -3.39%  '
$ws.Range('E38').Value = '  -2.47%  '
$ws.Range('D39').Value = '''0.525'
$ws.Range('E39').Value = '  -0.67%  '
$ws.Range('D40').Value = '''0.819'
$ws.Range('E40').Value = '  -0.28%  '
$ws.Range('E41').Value = '  +0.09%  '
$ws.Range('D42').Value = '''0.991'
$ws.Range('E42').Value = '  +1.37%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').Value = '''1.79'
$ws.Range('E43').Value = '  -0.23%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').Value = '''5.34'
$ws.Range('E44').Value = '  +1.78%  '
$ws.Range('E45').Value = '  +0.69%  '
$ws.Range('D46').Value = '''63.30'
$ws.Range('E46').Value = '  -1.01%  '
$ws.Range('D47').Value = '''1.702.13'
$ws.Range('E47').Value = '  +0.40%  '
$ws.Range('D48').Value = '''85.64'
$ws.Range('E48').Value = '  +0.25%  '
$ws.Range('D49').Value = '''0.0₇0991'
$ws.Range('E49').Value = '  -0.70%  '
$ws.Range('D50').Value = '''0.0952'
$ws.Range('E50').Value = '  -0.47%  '
$ws.Range('D51').Value = '''0.0491'
$ws.Range('E51').Value = '  -0.86%  '
